$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 1).Value = [double]"-3.7207827015474401E-3"
$ws.Cells.Item(2, 2).Value = [double]"3.7207827015459899E-3"
$ws.Cells.Item(2, 3).Value = [double]"0.48514992071525398"
$ws.Cells.Item(2, 4).Value = [double]"22.978005452784998"
$ws.Cells.Item(2, 5).Value = [double]"-22.9780054527942"
$ws.Cells.Item(2, 6).Value = [double]"-14.1228776736695"
$ws.Cells.Item(2, 7).Value = [double]"-4.1915215717182302E-3"
$ws.Cells.Item(2, 8).Value = [double]"4.1915215717164096E-3"
$ws.Cells.Item(2, 9).Value = [double]"0.16335079127382501"
$ws.Cells.Item(2, 10).Value = [double]"6.3675205592801101"
$ws.Cells.Item(2, 11).Value = [double]"-6.3675205592802602"
$ws.Cells.Item(2, 12).Value = [double]"-3.8979899187855702"
$ws.Cells.Item(3, 1).Value = [double]"-3.7207827015493002E-3"
$ws.Cells.Item(3, 2).Value = [double]"-3.7207827015468399E-3"
$ws.Cells.Item(3, 3).Value = [double]"0.48514992071525298"
$ws.Cells.Item(3, 4).Value = [double]"-22.978005452788299"
$ws.Cells.Item(3, 5).Value = [double]"-22.978005452799199"
$ws.Cells.Item(3, 6).Value = [double]"14.1228776736698"
$ws.Cells.Item(3, 7).Value = [double]"-4.19152157171986E-3"
$ws.Cells.Item(3, 8).Value = [double]"-4.1915215717171104E-3"
$ws.Cells.Item(3, 9).Value = [double]"0.16335079127382299"
$ws.Cells.Item(3, 10).Value = [double]"-6.3675205592801598"
$ws.Cells.Item(3, 11).Value = [double]"-6.3675205592803303"
$ws.Cells.Item(3, 12).Value = [double]"3.8979899187855902"
$ws.Cells.Item(4, 1).Value = [double]"3.72078270154835E-3"
$ws.Cells.Item(4, 2).Value = [double]"3.7207827015466898E-3"
$ws.Cells.Item(4, 3).Value = [double]"0.48514992071524099"
$ws.Cells.Item(4, 4).Value = [double]"22.978005452787801"
$ws.Cells.Item(4, 5).Value = [double]"22.978005452796399"
$ws.Cells.Item(4, 6).Value = [double]"14.1228776736697"
$ws.Cells.Item(4, 7).Value = [double]"4.1915215717190498E-3"
$ws.Cells.Item(4, 8).Value = [double]"4.1915215717170402E-3"
$ws.Cells.Item(4, 9).Value = [double]"0.16335079127382099"
$ws.Cells.Item(4, 10).Value = [double]"6.3675205592801296"
$ws.Cells.Item(4, 11).Value = [double]"6.3675205592802904"
$ws.Cells.Item(4, 12).Value = [double]"3.89798991878552"
$ws.Cells.Item(5, 1).Value = [double]"3.7207827015473599E-3"
$ws.Cells.Item(5, 2).Value = [double]"-3.7207827015454898E-3"
$ws.Cells.Item(5, 3).Value = [double]"0.48514992071524998"
$ws.Cells.Item(5, 4).Value = [double]"-22.978005452783901"
$ws.Cells.Item(5, 5).Value = [double]"22.978005452794601"
$ws.Cells.Item(5, 6).Value = [double]"-14.1228776736702"
$ws.Cells.Item(5, 7).Value = [double]"4.1915215717182302E-3"
$ws.Cells.Item(5, 8).Value = [double]"-4.1915215717160401E-3"
$ws.Cells.Item(5, 9).Value = [double]"0.16335079127382099"
$ws.Cells.Item(5, 10).Value = [double]"-6.3675205592800701"
$ws.Cells.Item(5, 11).Value = [double]"6.3675205592802602"
$ws.Cells.Item(5, 12).Value = [double]"-3.8979899187855498"

$ws.Range("A2:L5").Select()
